# Apply updated crypto market data to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.479.81"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.00"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.64"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2632"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06183"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07069"
$ws.Range("E10").Value = "  -1.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.665.43"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.84"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5892"
$ws.Range("E13").Value = "  -5.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.372"
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "74.97"
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.475.29"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006754"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.45"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.880.34"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.448"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.736"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.288"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.88"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.03"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.390"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.722"
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.95"
$ws.Range("E29").Value = "  +1.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.950"
$ws.Range("E30").Value = "  +4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07810"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.646"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9993"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04211"
$ws.Range("E34").Value = "  -7.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.601"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6096"
$ws.Range("E36").Value = "  +4.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9485"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.596"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8581"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.851"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01472"
$ws.Range("E42").Value = "  -6.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.84"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3767"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.837"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.205"
$ws.Range("E47").Value = "  +1.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05253"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "29.80"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.335"
$ws.Range("E51").Value = "  +1.51%  "
